# Auto-generated Excel COM-interop script applying the coin-price update diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric-looking cells (columns D and E) must keep their exact text
# representation (trailing zeros, leading zeros, percent signs, etc.),
# so force the cell to Text format before assigning the string value.
function Set-TextValue($ws, $cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

Set-TextValue $ws "D2" "303.75"
Set-TextValue $ws "E2" "0.61%"
Set-TextValue $ws "D3" "35.67"
Set-TextValue $ws "E3" "11.39%"
Set-TextValue $ws "D4" "5.092"
Set-TextValue $ws "E4" "1.08%"
Set-TextValue $ws "D5" "0.07801"
Set-TextValue $ws "E5" "0.24%"
Set-TextValue $ws "D6" "2.266"
Set-TextValue $ws "E6" "-2.06%"
Set-TextValue $ws "D7" "8.080"
Set-TextValue $ws "E7" "1.81%"
Set-TextValue $ws "D8" "4.051"
Set-TextValue $ws "E8" "5.16%"
Set-TextValue $ws "D9" "0.9282"
Set-TextValue $ws "E9" "-0.04%"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue $ws "D10" "0.1835"
Set-TextValue $ws "E10" "4.13%"
$ws.Range("B11").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C11").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue $ws "D11" "0.09001"
Set-TextValue $ws "E11" "-11.53%"
Set-TextValue $ws "D12" "0.08556"
Set-TextValue $ws "E12" "1.46%"
Set-TextValue $ws "D13" "0.03781"
Set-TextValue $ws "E13" "13.05%"
Set-TextValue $ws "D14" "0.09938"
Set-TextValue $ws "E14" "0.18%"
Set-TextValue $ws "D15" "0.001475"
Set-TextValue $ws "E15" "0.05%"
Set-TextValue $ws "D16" "0.005659"
Set-TextValue $ws "E16" "-1.89%"
Set-TextValue $ws "D17" "3.486"
Set-TextValue $ws "E17" "-0.29%"
Set-TextValue $ws "D18" "2.184"
Set-TextValue $ws "E18" "-0.33%"
Set-TextValue $ws "E19" "2.97%"
Set-TextValue $ws "D20" "0.1323"
Set-TextValue $ws "E20" "-1.57%"
Set-TextValue $ws "D21" "4.578"
Set-TextValue $ws "E21" "7.23%"
Set-TextValue $ws "D23" "0.04677"
Set-TextValue $ws "E23" "1.36%"
Set-TextValue $ws "D24" "0.001233"
Set-TextValue $ws "E24" "1.26%"
Set-TextValue $ws "D25" "0.004529"
Set-TextValue $ws "E25" "3.42%"
Set-TextValue $ws "E26" "1.09%"
Set-TextValue $ws "E27" "-20.03%"
Set-TextValue $ws "D39" "0.01767"
Set-TextValue $ws "E39" "3.32%"
Set-TextValue $ws "D40" "0.04738"
Set-TextValue $ws "E40" "0.07%"
Set-TextValue $ws "D41" "0.007980"
Set-TextValue $ws "E41" "2.05%"
Set-TextValue $ws "D42" "0.1414"
Set-TextValue $ws "E42" "0.68%"
Set-TextValue $ws "E43" "-18.20%"
Set-TextValue $ws "D44" "0.002305"
Set-TextValue $ws "E44" "11.78%"
Set-TextValue $ws "D45" "0.009622"
Set-TextValue $ws "E45" "-0.32%"
Set-TextValue $ws "D46" "0.00006202"
Set-TextValue $ws "E46" "1.59%"
Set-TextValue $ws "E47" "0.78%"
Set-TextValue $ws "D48" "5.803"
Set-TextValue $ws "E48" "118.62%"
Set-TextValue $ws "D49" "0.002697"
Set-TextValue $ws "E49" "35.55%"
Set-TextValue $ws "D50" "0.00002105"
Set-TextValue $ws "E50" "0.78%"
Set-TextValue $ws "D51" "0.0002005"
Set-TextValue $ws "E51" "0.78%"
